$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2654603333333334
$ws.Range("H2").Value = 0.796381
$ws.Range("I2").Value = 0.04079010536687974
$ws.Range("J2").Value = 0.04079010536687975
$ws.Range("M2").Value = 15.75563966666667
$ws.Range("N2").Value = 47.266919
$ws.Range("O2").Value = 0.3220556913988901
$ws.Range("P2").Value = 0.32205569139889
$ws.Range("Q2").Value = 4.182497357793222
$ws.Range("R2").Value = 37.642476220139
$ws.Range("S2").Value = 0.01313668558616403
$ws.Range("T2").Value = 0.01313668558616403
$ws.Range("G3").Value = 0.2654603333333334
$ws.Range("H3").Value = 0.796381
$ws.Range("I3").Value = 0.04079010536687974
$ws.Range("J3").Value = 0.04079010536687975
$ws.Range("O3").Value = 0.5509544596378365
$ws.Range("P3").Value = 0.5509544596378364
$ws.Range("Q3").Value = 7.155177297722445
$ws.Range("R3").Value = 64.396595679502
$ws.Range("S3").Value = 0.02247349046097964
$ws.Range("T3").Value = 0.02247349046097964
$ws.Range("G4").Value = 0.2654603333333334
$ws.Range("H4").Value = 0.796381
$ws.Range("I4").Value = 0.04079010536687974
$ws.Range("J4").Value = 0.04079010536687975
$ws.Range("O4").Value = 0.1269898489632735
$ws.Range("P4").Value = 0.1269898489632735
$ws.Range("Q4").Value = 1.649201432983222
$ws.Range("R4").Value = 14.842812896849
$ws.Range("S4").Value = 0.005179929319736072
$ws.Range("T4").Value = 0.005179929319736072
$ws.Range("I5").Value = 0.8420553458721338
$ws.Range("J5").Value = 0.8420553458721339
$ws.Range("M5").Value = 15.75563966666667
$ws.Range("N5").Value = 47.266919
$ws.Range("O5").Value = 0.3220556913988901
$ws.Range("P5").Value = 0.32205569139889
$ws.Range("Q5").Value = 86.34187697111278
$ws.Range("R5").Value = 777.076892740015
$ws.Range("S5").Value = 0.2711887166109815
$ws.Range("T5").Value = 0.2711887166109815
$ws.Range("I6").Value = 0.8420553458721338
$ws.Range("J6").Value = 0.8420553458721339
$ws.Range("O6").Value = 0.5509544596378365
$ws.Range("P6").Value = 0.5509544596378364
$ws.Range("S6").Value = 0.463934148070133
$ws.Range("T6").Value = 0.463934148070133
$ws.Range("I7").Value = 0.8420553458721338
$ws.Range("J7").Value = 0.8420553458721339
$ws.Range("O7").Value = 0.1269898489632735
$ws.Range("P7").Value = 0.1269898489632735
$ws.Range("S7").Value = 0.1069324811910193
$ws.Range("T7").Value = 0.1069324811910193
$ws.Range("G8").Value = 0.7624369999999999
$ws.Range("I8").Value = 0.1171545487609863
$ws.Range("J8").Value = 0.1171545487609864
$ws.Range("M8").Value = 15.75563966666667
$ws.Range("N8").Value = 47.266919
$ws.Range("O8").Value = 0.3220556913988901
$ws.Range("P8").Value = 0.32205569139889
$ws.Range("Q8").Value = 12.01268264053433
$ws.Range("R8").Value = 108.114143764809
$ws.Range("S8").Value = 0.03773028920174443
$ws.Range("T8").Value = 0.03773028920174443
$ws.Range("G9").Value = 0.7624369999999999
$ws.Range("I9").Value = 0.1171545487609863
$ws.Range("J9").Value = 0.1171545487609864
$ws.Range("O9").Value = 0.5509544596378365
$ws.Range("P9").Value = 0.5509544596378364
$ws.Range("Q9").Value = 20.55061049928466
$ws.Range("S9").Value = 0.0645468211067238
$ws.Range("T9").Value = 0.06454682110672379
$ws.Range("G10").Value = 0.7624369999999999
$ws.Range("I10").Value = 0.1171545487609863
$ws.Range("J10").Value = 0.1171545487609864
$ws.Range("O10").Value = 0.1269898489632735
$ws.Range("P10").Value = 0.1269898489632735
$ws.Range("Q10").Value = 4.736723476424332
$ws.Range("R10").Value = 14.842812896849
$ws.Range("S10").Value = 0.01487743845251812
$ws.Range("T10").Value = 0.01487743845251812
